$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 75 was blank; fill in a completed time-log entry (Coding, 10/6/2014).
$ws.Range("A75").Value = 41918
$ws.Range("B75").Value = 0.59375
$ws.Range("C75").Value = 0.61736111111111114
$ws.Range("D75").Value = 5
$ws.Range("F75").Value = "Coding"

$ws.Range("A74:D74").Copy() | Out-Null
$ws.Range("A75:D75").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("F74").Copy() | Out-Null
$ws.Range("F75").PasteSpecial(-4122) | Out-Null        # xlPasteFormats
$excel.CutCopyMode = $false

# Force the dependent shared-formula cell to recompute; it otherwise keeps
# a stale cached value from before D75 had a value.
$ws.Range("E75").Formula = $ws.Range("E75").Formula

$ws.Range("A76").Select() | Out-Null
